$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text assignments (safe: not ambiguous with numeric parsing)
$ws.Range('D2').Value = '39.374.36'
$ws.Range('E2').Value = '  -2.20%  '
$ws.Range('D3').Value = '2.203.08'
$ws.Range('E3').Value = '  -6.02%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('E5').Value = '  -4.27%  '
$ws.Range('E6').Value = '  -4.71%  '
$ws.Range('E7').Value = '  -3.67%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -3.66%  '
$ws.Range('E10').Value = '  -6.23%  '
$ws.Range('E11').Value = '  -4.39%  '
$ws.Range('E12').Value = '  -11.00%  '
$ws.Range('E13').Value = '  -2.75%  '
$ws.Range('D14').Value = '2.539.78'
$ws.Range('E14').Value = '  -6.23%  '
$ws.Range('E15').Value = '  -3.85%  '
$ws.Range('E16').Value = '  -5.90%  '
$ws.Range('D17').Value = '2.188.89'
$ws.Range('E17').Value = '  -6.81%  '
$ws.Range('E18').Value = '  -6.17%  '
$ws.Range('D19').Value = '39.264.00'
$ws.Range('E19').Value = '  -2.40%  '
$ws.Range('D20').Value = '0.0₃0869'
$ws.Range('E20').Value = '  -4.19%  '
$ws.Range('E21').Value = '  -6.52%  '
$ws.Range('E22').Value = '  -4.49%  '
$ws.Range('E23').Value = '  -4.63%  '
$ws.Range('E24').Value = '  -4.77%  '
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('E26').Value = '  -6.55%  '
$ws.Range('E27').Value = '  -0.82%  '
$ws.Range('E28').Value = '  -4.26%  '
$ws.Range('E29').Value = '  +0.77%  '
$ws.Range('E30').Value = '  -2.08%  '
$ws.Range('E31').Value = '  -2.27%  '
$ws.Range('E32').Value = '  -9.49%  '
$ws.Range('E33').Value = '  -0.19%  '
$ws.Range('E34').Value = '  -7.00%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('E35').Value = '  -4.69%  '
$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('E36').Value = '  -4.81%  '
$ws.Range('E37').Value = '  -3.27%  '
$ws.Range('E38').Value = '  -3.47%  '
$ws.Range('E39').Value = '  -5.79%  '
$ws.Range('E40').Value = '  -4.74%  '
$ws.Range('E41').Value = '  -4.37%  '
$ws.Range('E42').Value = '  -5.41%  '
$ws.Range('D43').Value = '1.901.68'
$ws.Range('E43').Value = '  -2.70%  '
$ws.Range('E44').Value = '  -9.19%  '
$ws.Range('E45').Value = '  -3.46%  '
$ws.Range('E46').Value = '  -3.86%  '
$ws.Range('E47').Value = '  -9.79%  '
$ws.Range('E48').Value = '  -4.02%  '
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('D50').Value = '2.406.70'
$ws.Range('E50').Value = '  -6.38%  '
$ws.Range('E51').Value = '  -6.38%  '

# Numeric-looking text values: use formula-literal + paste-as-values so Excel
# keeps them as text (matching original inlineStr cells) without adding any
# NumberFormat/quotePrefix style to the cell.
$ws.Range('D5').Formula = '="295.50"'
$ws.Range('D6').Formula = '="81.77"'
$ws.Range('D7').Formula = '="0.511"'
$ws.Range('D9').Formula = '="0.469"'
$ws.Range('D10').Formula = '="0.0770"'
$ws.Range('D11').Formula = '="29.12"'
$ws.Range('D12').Formula = '="47.04"'
$ws.Range('D15').Formula = '="6.23"'
$ws.Range('D16').Formula = '="13.94"'
$ws.Range('D18').Formula = '="0.711"'
$ws.Range('D22').Formula = '="64.82"'
$ws.Range('D23').Formula = '="10.30"'
$ws.Range('D26').Formula = '="2.40"'
$ws.Range('D27').Formula = '="1.80"'
$ws.Range('D28').Formula = '="22.56"'
$ws.Range('D30').Formula = '="9.07"'
$ws.Range('D31').Formula = '="148.75"'
$ws.Range('D32').Formula = '="31.78"'
$ws.Range('D33').Formula = '="0.999"'
$ws.Range('D34').Formula = '="4.79"'
$ws.Range('D35').Formula = '="0.0695"'
$ws.Range('D36').Formula = '="2.33"'
$ws.Range('D38').Formula = '="15.36"'
$ws.Range('D40').Formula = '="2.63"'
$ws.Range('D41').Formula = '="1.65"'
$ws.Range('D42').Formula = '="3.62"'
$ws.Range('D45').Formula = '="0.0259"'
$ws.Range('D46').Formula = '="8.98"'
$ws.Range('D49').Formula = '="71.81"'
$ws.Range('D51').Formula = '="87.28"'

$rng = $ws.Range("D5","D6","D7","D9","D10","D11","D12","D15","D16","D18","D22","D23","D26","D27","D28","D30","D31","D32","D33","D34","D35","D36","D38","D40","D41","D42","D45","D46","D49","D51")
$rng.Copy()
$rng.PasteSpecial(-4163)
$excel.CutCopyMode = 0
